# Updated cryptos list values (price + 1h volume change) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.036.50'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '2.007.59'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.90'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.373'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.101'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.07%  '
$ws.Range('D12').Value = '2.305.53'
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.733'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('D17').Value = '2.016.21'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '36.955.84'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').Value = '0.0₃0809'
$ws.Range('E21').Value = '  -3.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.123'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.116'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0596'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.31%  '
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').Value = '1.451.41'
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '94.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.90%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0905'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.32%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.43%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.32%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.992'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('D51').Value = '2.196.43'
